$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.145.12"

$ws.Range("D3").Value = "3.521.61"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'607.86"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'148.09"
$ws.Range("E6").Value = "  -2.64%  "

$ws.Range("D7").Value = "3.520.75"
$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -1.83%  "

$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").Value = "'7.87"
$ws.Range("E11").Value = "  +2.71%  "

$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").Value = "4.119.24"
$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("D15").Value = "'31.89"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "3.518.81"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "67.110.61"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "'10.74"
$ws.Range("E19").Value = "  +9.04%  "

$ws.Range("D20").Value = "'6.43"
$ws.Range("E20").Value = "  -1.29%  "

$ws.Range("D21").Value = "'15.36"
$ws.Range("E21").Value = "  -0.81%  "

$ws.Range("D22").Value = "'438.21"
$ws.Range("E22").Value = "  -2.01%  "

$ws.Range("D23").Value = "'0.610"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").Value = "'79.60"
$ws.Range("E24").Value = "  +1.55%  "

$ws.Range("D25").Value = "3.670.66"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -4.78%  "

$ws.Range("D28").Value = "'9.80"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").Value = "'8.29"
$ws.Range("E29").Value = "  -5.60%  "

$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("E31").Value = "  -3.17%  "

$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = "  -2.55%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "'25.45"
$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("D35").Value = "3.518.14"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("D37").Value = "'5.94"
$ws.Range("E37").Value = "  -3.36%  "

$ws.Range("D38").Value = "'8.06"
$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").Value = "'0.0895"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").Value = "'172.24"
$ws.Range("E42").Value = "  -3.56%  "

$ws.Range("D43").Value = "'5.45"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "'2.08"
$ws.Range("E44").Value = "  -9.84%  "

$ws.Range("D45").Value = "'0.896"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").Value = "'46.05"
$ws.Range("E46").Value = "  -0.90%  "

$ws.Range("D47").Value = "'28.59"
$ws.Range("E47").Value = "  -5.45%  "

$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("D49").Value = "'7.49"
$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("D50").Value = "'2.45"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("E51").Value = "  +0.33%  "
